# The document opened with an extra, empty, formatted paragraph above the
# actual title paragraph. Remove that leading empty paragraph so the
# "Psicología infantil..." paragraph becomes the first paragraph in the body.
$d = $word.ActiveDocument

$firstPara = $d.Paragraphs.Item(1)
if ($firstPara.Range.Text.Trim().Length -eq 0) {
    $firstPara.Range.Delete()
}

# Hide the built-in "Default Paragraph Font" character style from the
# style gallery / recommended list (adds w:semiHidden to its definition).
$defParaFontStyle = $d.Styles.Item("Default Paragraph Font")
$defParaFontStyle.Visibility = $false
